$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table with refreshed values from the feed.
# Price cells (column D) that look like plain numbers are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# source data, e.g. "1.00", "0.0000322") instead of silently coercing them
# to numbers (which would drop formatting like trailing zeros or multiple
# dot separators used as thousands markers, e.g. "71.439.37").

# Row 2
$ws.Range('D2').Value = '71.439.37'
$ws.Range('E2').Value = '  -1.26%  '
# Row 3
$ws.Range('D3').Value = '3.966.27'
$ws.Range('E3').Value = '  -2.78%  '
# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.04%  '
# Row 5
$ws.Range('D5').Value = '''538.82'
$ws.Range('E5').Value = '  +3.13%  '
# Row 6
$ws.Range('D6').Value = '''149.96'
$ws.Range('E6').Value = '  +0.85%  '
# Row 7
$ws.Range('D7').Value = '3.960.35'
$ws.Range('E7').Value = '  -2.65%  '
# Row 8
$ws.Range('D8').Value = '''0.689'
# Row 9
$ws.Range('E9').Value = '  -0.11%  '
# Row 10
$ws.Range('D10').Value = '''0.747'
$ws.Range('E10').Value = '  -4.19%  '
# Row 11
$ws.Range('D11').Value = '''0.168'
$ws.Range('E11').Value = '  -6.67%  '
# Row 12
$ws.Range('D12').Value = '''55.25'
$ws.Range('E12').Value = '  +13.18%  '
# Row 13
$ws.Range('D13').Value = '''0.0000322'
$ws.Range('E13').Value = '  -3.88%  '
# Row 14
$ws.Range('D14').Value = '''10.75'
$ws.Range('E14').Value = '  -2.79%  '
# Row 15
$ws.Range('D15').Value = '4.590.32'
$ws.Range('E15').Value = '  -2.77%  '
# Row 16
$ws.Range('D16').Value = '3.957.26'
$ws.Range('E16').Value = '  -2.32%  '
# Row 17
$ws.Range('D17').Value = '''14.05'
$ws.Range('E17').Value = '  -3.48%  '
# Row 18
$ws.Range('E18').Value = '  -3.84%  '
# Row 19
$ws.Range('E19').Value = '  -1.75%  '
# Row 20
$ws.Range('D20').Value = '''1.18'
$ws.Range('E20').Value = '  -5.72%  '
# Row 21
$ws.Range('D21').Value = '71.253.62'
$ws.Range('E21').Value = '  -1.51%  '
# Row 22
$ws.Range('D22').Value = '''430.11'
$ws.Range('E22').Value = '  -4.17%  '
# Row 23
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = '''3.60'
$ws.Range('E23').Value = '  -0.72%  '
# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '''97.57'
$ws.Range('E24').Value = '  -6.32%  '
# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '''4.26'
$ws.Range('E25').Value = '  +4.71%  '
# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '''14.67'
$ws.Range('E26').Value = '  -2.55%  '
# Row 27
$ws.Range('D27').Value = '''11.45'
$ws.Range('E27').Value = '  -0.02%  '
# Row 28
$ws.Range('D28').Value = '''3.97'
$ws.Range('E28').Value = '  +20.35%  '
# Row 29
$ws.Range('D29').Value = '''10.82'
$ws.Range('E29').Value = '  -2.95%  '
# Row 30
$ws.Range('E30').Value = '  +1.26%  '
# Row 31
$ws.Range('D31').Value = '''36.82'
$ws.Range('E31').Value = '  -3.43%  '
# Row 32
$ws.Range('D32').Value = '''7.97'
$ws.Range('E32').Value = '  +19.52%  '
# Row 33
$ws.Range('D33').Value = '''50.97'
$ws.Range('E33').Value = '  +19.90%  '
# Row 34
$ws.Range('D34').Value = '''0.133'
$ws.Range('E34').Value = '  +0.29%  '
# Row 35
$ws.Range('D35').Value = '''13.44'
$ws.Range('E35').Value = '  -2.19%  '
# Row 36
$ws.Range('D36').Value = '''682.47'
$ws.Range('E36').Value = '  -0.18%  '
# Row 37
$ws.Range('D37').Value = '''65.73'
$ws.Range('E37').Value = '  -3.10%  '
# Row 38
$ws.Range('D38').Value = '''0.443'
$ws.Range('E38').Value = '  +1.97%  '
# Row 39
$ws.Range('D39').Value = '0.0₃0826'
$ws.Range('E39').Value = '  -7.68%  '
# Row 40
$ws.Range('E40').Value = '  -3.41%  '
# Row 41
$ws.Range('D41').Value = '''3.42'
$ws.Range('E41').Value = '  -1.38%  '
# Row 42
$ws.Range('E42').Value = '  +0.12%  '
# Row 43
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.07%  '
# Row 44
$ws.Range('D44').Value = '''0.0486'
$ws.Range('E44').Value = '  -3.80%  '
# Row 45
$ws.Range('D45').Value = '''3.21'
$ws.Range('E45').Value = '  -0.61%  '
# Row 46
$ws.Range('D46').Value = '''10.33'
$ws.Range('E46').Value = '  +3.97%  '
# Row 47
$ws.Range('E47').Value = '  -5.33%  '
# Row 48
$ws.Range('D48').Value = '''2.67'
$ws.Range('E48').Value = '  -1.04%  '
# Row 49
$ws.Range('D49').Value = '''3.36'
$ws.Range('E49').Value = '  -1.74%  '
# Row 50
$ws.Range('D50').Value = '''3.02'
$ws.Range('E50').Value = '  -2.28%  '
# Row 51
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '''0.000269'
$ws.Range('E51').Value = '  -6.00%  '
